$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("M2").Value = 7.655977
$ws.Range("N2").Value = 22.967931
$ws.Range("O2").Value = 0.2994795900616967
$ws.Range("P2").Value = 0.2994795900616967
$ws.Range("Q2").Value = 34.981468085067
$ws.Range("R2").Value = 314.833212765603
$ws.Range("S2").Value = 0.2994795900616967
$ws.Range("T2").Value = 0.2994795900616967

# Row 3 updates
$ws.Range("O3").Value = 0.3140620915319453
$ws.Range("P3").Value = 0.3140620915319453
$ws.Range("S3").Value = 0.3140620915319453
$ws.Range("T3").Value = 0.3140620915319453

# Row 4 updates
$ws.Range("M4").Value = 9.879524666666667
$ws.Range("N4").Value = 29.638574
$ws.Range("O4").Value = 0.386458318406358
$ws.Range("P4").Value = 0.386458318406358
$ws.Range("Q4").Value = 45.141237600718
$ws.Range("R4").Value = 406.271138406462
$ws.Range("S4").Value = 0.386458318406358
$ws.Range("T4").Value = 0.386458318406358

$wb.Save()
